# Update cryptocurrency price and 1h-volume-change data on the active
# worksheet of the active workbook (scheduled GitHub Actions refresh of
# cryptos.xlsx). Price cells are forced to Text format before/while being
# written so values such as "581.65" are stored as text (matching the
# original inline-string data) instead of being auto-converted to numbers,
# then the style is reset to "Normal" so no stray number-format/style is
# left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.566.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.59%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.558.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.10%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.97%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.629"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.56%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.546.08"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.93%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("E10").Value = "  +20.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.650"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.63%  "
$ws.Range("E13").Value = "  +5.95%  "
$ws.Range("E14").Value = "  +0.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.126.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "70.698.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.570.68"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "570.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.06%  "
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("E22").Value = "  -0.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.86%  "
$ws.Range("E24").Value = "  +3.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "94.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.26%  "
$ws.Range("E28").Value = "  +2.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.58%  "
$ws.Range("E32").Value = "  -1.75%  "
$ws.Range("E33").Value = "  +2.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "62.96"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.82%  "
$ws.Range("E35").Value = "  +13.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "549.88"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.80%  "
$ws.Range("E37").Value = "  +5.05%  "
$ws.Range("E38").Value = "  +10.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.80"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.37%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0798"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.601.51"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +12.27%  "
$ws.Range("E43").Value = "  +3.69%  "
$ws.Range("E44").Value = "  +3.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0463"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.67%  "
$ws.Range("E46").Value = "  +1.41%  "
$ws.Range("E47").Value = "  -1.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.30%  "
$ws.Range("E49").Value = "  +3.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +14.77%  "
$ws.Range("E51").Value = "  +0.00%  "

